# Incorporate updated data from upstream processes through 2024.
#
# The workbook's "Sheet1" tracks Delaware facility counts by technology
# (columns B-G) for each "Open year" (column A, row 2 = 2000 ... row 26 = 2024).
# The chart on the sheet reads its series directly from these columns, so
# updating the cells is sufficient to reflect the new figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 2023 row (row 25): Solar (column E) 87 -> 88
$ws.Range("E25").Value = 88

# 2024 row (row 26): Energy Storage (column C) 1 -> 2, Solar (column E) 47 -> 78
$ws.Range("C26").Value = 2
$ws.Range("E26").Value = 78
